# Circle Language Spec: System Objects: Assignment
# Marks common "trigger" words (is/are/same/most/clear/always/requir) in
# red throughout the document, and relocates the _GoBack bookmark to sit
# inside the first edit ("will" -> "wi" + bookmark + "ll"), matching the
# way Word stamps the last-edit-location bookmark at save time.

$d = $word.ActiveDocument
$wdColorRed = 255

# ---------------------------------------------------------------------
# Helper: re-key the full text of a (unique) anchor phrase into a single
# run, so that any pre-existing run split inside the target sub-word
# (an artifact of earlier edits) does not fragment the new colored run.
# ---------------------------------------------------------------------
function Normalize-Anchor($anchor) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NORMALIZE: NOT FOUND: $anchor"
        return
    }
    $rng.Text = "X"
    $rng2 = $d.Range($rng.Start, $rng.Start + 1)
    $rng2.Text = $anchor
}

# ---------------------------------------------------------------------
# Helper: find a unique anchor phrase, locate a sub-word inside it
# (whole-word match unless $plain is set, for partial-word splits like
# "requir"/"e"), and color that sub-range red.
# ---------------------------------------------------------------------
function Color-SubWord($anchor, $word, $plain) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "COLOR: NOT FOUND: $anchor"
        return
    }
    $base = $rng.Start
    if ($plain) {
        $relOffset = $anchor.IndexOf($word)
    } else {
        $m = [regex]::Match($anchor, "\b" + [regex]::Escape($word) + "\b")
        if (-not $m.Success) {
            Write-Output "COLOR: WORD NOT FOUND: $word in $anchor"
            return
        }
        $relOffset = $m.Index
    }
    if ($relOffset -lt 0) {
        Write-Output "COLOR: WORD NOT FOUND (plain): $word in $anchor"
        return
    }
    $wordStart = $base + $relOffset
    $wordEnd = $wordStart + $word.Length
    $wordRng = $d.Range($wordStart, $wordEnd)
    $wordRng.Font.Color = $wdColorRed
}

# 1. "You will usually not see " -> "You " + red("wi") + _GoBack + red("ll") + " usually..."
$rng = $d.Content
$found = $rng.Find.Execute("will", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $wiStart = $rng.Start
    $wiEnd = $rng.Start + 2
    $llStart = $wiEnd
    $llEnd = $rng.End

    $wiRng = $d.Range($wiStart, $wiEnd)
    $wiRng.Font.Color = $wdColorRed
    $llRng = $d.Range($llStart, $llEnd)
    $llRng.Font.Color = $wdColorRed

    $bmRng = $d.Range($wiEnd, $wiEnd)
    $d.Bookmarks.Add("_GoBack", $bmRng)
} else {
    Write-Output "OP1: 'will' not found"
}

# 2. "commands. Those system commands are called " -> color "are"
Color-SubWord "commands. Those system commands are called " "are" $false

# 3. "Different aspects ... Below is an overview of the most common types of assignments."
$anchor3 = "Different aspects have different types of assignment. Below is an overview of the most common types of assignments."
Normalize-Anchor $anchor3
Color-SubWord $anchor3 "is" $false
Color-SubWord $anchor3 "most" $false

# 4. "It is also made clear in the overview, which " -> color "is" and "clear"
$anchor4 = "It is also made clear in the overview, which "
Color-SubWord $anchor4 "is" $false
Color-SubWord $anchor4 "clear" $false

# 5. "commands are called to perform the assignment." -> color "are"
Color-SubWord "commands are called to perform the assignment." "are" $false

# 6. "Object-bound aspects ... are displayed differently. When a reference-bound aspect is "
$anchor6 = "Object-bound aspects and reference-bound aspects are displayed differently. When a reference-bound aspect is "
Color-SubWord $anchor6 "are" $false
Color-SubWord $anchor6 "is" $false

# 7. " then the reference is displayed with a parent around it:" -> color "is"
Color-SubWord " then the reference is displayed with a parent around it:" "is" $false

# 8. "When an object-bound aspect is " -> color "is"
Color-SubWord "When an object-bound aspect is " "is" $false

# 9. " then the targeted object is displayed without a parent around it:" -> color "is"
Color-SubWord " then the targeted object is displayed without a parent around it:" "is" $false

# 10. "the same object as the source." (italic run) -> color "same"
Color-SubWord "the same object as the source." "same" $false

# 11. "In the assignment notation the line type indicates which aspect is " -> color "is"
Color-SubWord "In the assignment notation the line type indicates which aspect is " "is" $false

# 12. "assignment does not require an assignment call symbol, because a " -> color "requir"
Color-SubWord "assignment does not require an assignment call symbol, because a " "requir" $true

# 13. "connection is always an assignment." -> color "is" and "always"
$anchor13 = "connection is always an assignment."
Color-SubWord $anchor13 "is" $false
Color-SubWord $anchor13 "always" $false

# 14. "aspect in an assignment is to get the " -> color "is"
Color-SubWord "aspect in an assignment is to get the " "is" $false

# 15. ". Less conventional ways of yielding over aspects like that, is also called a " -> color "is"
Color-SubWord ". Less conventional ways of yielding over aspects like that, is also called a " "is" $false

# 16. "Also note here, that there are two ways to get the " -> color "are"
Color-SubWord "Also note here, that there are two ways to get the " "are" $false

# 17. "If the source of the assignment is a pointer-to-pointer, then the target also becomes a pointer-to-pointer. So this also gives " -> color "is"
Color-SubWord "If the source of the assignment is a pointer-to-pointer, then the target also becomes a pointer-to-pointer. So this also gives " "is" $false

# 18. Remove the old trailing _GoBack bookmark (its position has moved to
#     inside "will", handled in step 1 above).
$exists = $d.Bookmarks.Exists("_GoBack")
if ($exists) {
    # Only remove it if it is not the one we just added inside "will".
    $bm = $d.Bookmarks.Item("_GoBack")
    if ($bm.Start -ne $wiEnd) {
        $bm.Delete()
    }
}

Write-Output "Done"
